$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'60.570.50"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -3.31%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'3.344.50"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -2.94%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("D5").Value = "'566.90"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -2.20%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'146.46"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -1.26%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "'  +0.04%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = "'  +0.36%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'7.91"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -1.18%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("E10").Value = "'  -1.48%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.414"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +1.15%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'3.913.83"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -3.05%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("E13").Value = "'  +1.03%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'27.69"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -2.21%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'3.334.50"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -3.39%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("E16").Value = "'  -1.91%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'60.602.24"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -3.33%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'6.27"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -0.94%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'14.54"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -0.79%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'8.89"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -1.92%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'376.87"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -2.63%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("E22").Value = "'  -0.86%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'74.73"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -0.86%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = "'  +0.00%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("E26").Value = "'  -5.77%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = "'  -4.61%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("E28").Value = "'  +0.12%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'7.33"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -4.29%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("E30").Value = "'  -1.54%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D32").Value = "'7.68"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -4.06%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'22.89"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  -1.53%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'1.29"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -3.59%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'5.29"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -1.57%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'1.55"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -5.10%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("B37").Value = "'Monero"
$ws.Range("B37").Style = "Normal"
$ws.Range("C37").Value = "'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("C37").Style = "Normal"
$ws.Range("D37").Value = "'167.39"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -1.25%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("B38").Value = "'Aptos"
$ws.Range("B38").Style = "Normal"
$ws.Range("C38").Value = "'https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("C38").Style = "Normal"
$ws.Range("D38").Value = "'6.80"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -2.06%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'27.95"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -12.36%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'3.379.48"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -2.98%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.0746"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -3.66%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("E42").Value = "'  -3.68%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("E43").Value = "'  -1.77%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("E44").Value = "'  -3.24%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("E45").Value = "'  -5.04%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'2.454.48"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -4.56%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = "'  -3.40%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("E48").Value = "'  +0.05%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'22.32"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -1.70%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Value = "'  -2.11%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'0.815"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -0.07%  "
$ws.Range("E51").Style = "Normal"
